$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Range("D1").Value = "BenchPress"
$ws.Activate()
$ws.Range("D1").Select()
